{"js": "// The first paragraph in the document body is a single, green-colored run\n// that currently holds the stale \"[LLM error: ...]\" placeholder text. Replace\n// its text with the commit summary, using manual line breaks (represented as\n// \"\\u000B\" in Office.js text strings, which Word serializes as <w:br/>)\n// between the summary lines so everything stays inside that one run.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\nconst lineBreak = \"\\u000B\";\nconst newText =\n  \"The main changes between the two versions of the document are:\" +\n  lineBreak + lineBreak +\n  \"- Two new test cases (TC-new and TC099999999-new) were added to validate passwords with less than 8 characters.\" +\n  lineBreak + lineBreak +\n  \"- A new section 3 was added with additional password inputs to test based on the new special character requirement from the client.\";\n\ntargetParagraph.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The target paragraph is the very first paragraph in the document body,\n# a single run colored green that currently holds the old \"[LLM error: ...]\" text.\n$p = $d.Paragraphs(1)\n$r = $p.Range\n\n# Word represents a manual line break (<w:br/>) inside a run's text as the\n# vertical-tab character (Chr(11)) when read/written through Range.Text.\n$vt = [char]11\n\n$newText = \"The main changes between the two versions of the document are:\" + $vt + $vt + `\n    \"- Two new test cases (TC-new and TC099999999-new) were added to validate passwords with less than 8 characters.\" + $vt + $vt + `\n    \"- A new section 3 was added with additional password inputs to test based on the new special character requirement from the client.\"\n\n$r.Text = $newText\n"}
